$d = $word.ActiveDocument
$t = $d.Tables(1)

# Map of 1-based row index -> new cell(1) text for rows whose content changed.
$changes = @{
    1  = "0M"
    2  = "0M"
    3  = "0M"
    4  = "1504"
    5  = "0.00001"
    6  = "0.00058"
    7  = "0.00012"
    9  = "0.00017"
    10 = "0.00019"
    11 = "0.00022"
    12 = "0.19321"
    44 = "99.93"
    45 = "0.19"
    46 = "288"
}

foreach ($rowIndex in $changes.Keys) {
    $t.Rows($rowIndex).Cells(1).Range.Text = $changes[$rowIndex]
}
